$wb = $excel.ActiveWorkbook

# Update "Correspond Handoff Datetime" / "Correspond Handback DateTime" values
# on the zh-cn worksheet for the dda6627d... handoff row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-18 04:20:47"
$wsZhCn.Range("H4").Value = "2016-03-18 04:21:10"

# Same update on the de-de worksheet.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-18 04:20:49"
$wsDeDe.Range("H4").Value = "2016-03-18 04:21:15"
